$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows before row 2; this shifts the existing 20 data rows
# (originally rows 2-21, timestamps 0-1900) down to rows 11-30.
$ws.Rows("2:10").Insert()
$ws.Rows("2:10").ClearFormats()

# Populate the 9 newly inserted rows with new sensor readings (timestamps 0-800)
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "falling"
$ws.Cells.Item(2, 3).Value = -0.4484343528747558
$ws.Cells.Item(2, 4).Value = 0.5805449485778809
$ws.Cells.Item(2, 5).Value = -0.5587977170944214
$ws.Cells.Item(2, 6).Value = 0.0189368221908807
$ws.Cells.Item(2, 7).Value = -0.030695978552103
$ws.Cells.Item(2, 8).Value = -0.0142026171088218

$ws.Cells.Item(3, 1).Value = 100
$ws.Cells.Item(3, 2).Value = "falling"
$ws.Cells.Item(3, 3).Value = -0.6014323234558105
$ws.Cells.Item(3, 4).Value = 0.5997557640075684
$ws.Cells.Item(3, 5).Value = -0.4316743612289428
$ws.Cells.Item(3, 6).Value = -0.0242818929255008
$ws.Cells.Item(3, 7).Value = 0.0024434609804302
$ws.Cells.Item(3, 8).Value = 0.0209221355617046

$ws.Cells.Item(4, 1).Value = 200
$ws.Cells.Item(4, 2).Value = "falling"
$ws.Cells.Item(4, 3).Value = -0.2612781524658203
$ws.Cells.Item(4, 4).Value = 0.5193090438842773
$ws.Cells.Item(4, 5).Value = -0.4944255352020263
$ws.Cells.Item(4, 6).Value = -0.0287106670439243
$ws.Cells.Item(4, 7).Value = -0.0282525178045034
$ws.Cells.Item(4, 8).Value = 0.00534507073462

$ws.Cells.Item(5, 1).Value = 300
$ws.Cells.Item(5, 2).Value = "falling"
$ws.Cells.Item(5, 3).Value = -0.2169137001037597
$ws.Cells.Item(5, 4).Value = 0.3676133155822754
$ws.Cells.Item(5, 5).Value = -0.733814001083374
$ws.Cells.Item(5, 6).Value = 0.0103847095742821
$ws.Cells.Item(5, 7).Value = -0.0001527163112768
$ws.Cells.Item(5, 8).Value = 0.0277943685650825

$ws.Cells.Item(6, 1).Value = 400
$ws.Cells.Item(6, 2).Value = "falling"
$ws.Cells.Item(6, 3).Value = -0.3025293350219726
$ws.Cells.Item(6, 4).Value = 0.4683008193969726
$ws.Cells.Item(6, 5).Value = -0.587003767490387
$ws.Cells.Item(6, 6).Value = -0.0300851128995418
$ws.Cells.Item(6, 7).Value = -0.0401643887162208
$ws.Cells.Item(6, 8).Value = 0.0236710291355848

$ws.Cells.Item(7, 1).Value = 500
$ws.Cells.Item(7, 2).Value = "falling"
$ws.Cells.Item(7, 3).Value = -0.4149298667907715
$ws.Cells.Item(7, 4).Value = 0.4677276611328125
$ws.Cells.Item(7, 5).Value = -0.731619656085968
$ws.Cells.Item(7, 6).Value = -0.0271835029125213
$ws.Cells.Item(7, 7).Value = -0.0343611687421798
$ws.Cells.Item(7, 8).Value = 0.0047342055477201

$ws.Cells.Item(8, 1).Value = 600
$ws.Cells.Item(8, 2).Value = "falling"
$ws.Cells.Item(8, 3).Value = -0.5511326789855957
$ws.Cells.Item(8, 4).Value = 0.6498098373413086
$ws.Cells.Item(8, 5).Value = -0.522668182849884
$ws.Cells.Item(8, 6).Value = -0.0245873257517814
$ws.Cells.Item(8, 7).Value = -0.0226020142436027
$ws.Cells.Item(8, 8).Value = -0.0163406450301408

$ws.Cells.Item(9, 1).Value = 700
$ws.Cells.Item(9, 2).Value = "falling"
$ws.Cells.Item(9, 3).Value = -0.5137066841125488
$ws.Cells.Item(9, 4).Value = 0.4998054504394531
$ws.Cells.Item(9, 5).Value = -0.6402766704559326
$ws.Cells.Item(9, 6).Value = -0.0479529201984405
$ws.Cells.Item(9, 7).Value = 0.07635815441608421
$ws.Cells.Item(9, 8).Value = -0.1252273768186569

$ws.Cells.Item(10, 1).Value = 800
$ws.Cells.Item(10, 2).Value = "falling"
$ws.Cells.Item(10, 3).Value = -0.6838326454162598
$ws.Cells.Item(10, 4).Value = 0.6059346199035645
$ws.Cells.Item(10, 5).Value = -0.2089821100234985
$ws.Cells.Item(10, 6).Value = 0.1218676194548606
$ws.Cells.Item(10, 7).Value = 0.3381139039993286
$ws.Cells.Item(10, 8).Value = -0.0650571510195732

# The shifted (originally existing) rows 11-30 keep their C:H readings untouched by
# the Insert above, but their timestamp (column A) needs to continue the sequence
# rather than restart at 0, so update column A for those rows to 900-2800.
$ws.Cells.Item(11, 1).Value = 900
$ws.Cells.Item(12, 1).Value = 1000
$ws.Cells.Item(13, 1).Value = 1100
$ws.Cells.Item(14, 1).Value = 1200
$ws.Cells.Item(15, 1).Value = 1300
$ws.Cells.Item(16, 1).Value = 1400
$ws.Cells.Item(17, 1).Value = 1500
$ws.Cells.Item(18, 1).Value = 1600
$ws.Cells.Item(19, 1).Value = 1700
$ws.Cells.Item(20, 1).Value = 1800
$ws.Cells.Item(21, 1).Value = 1900
$ws.Cells.Item(22, 1).Value = 2000
$ws.Cells.Item(23, 1).Value = 2100
$ws.Cells.Item(24, 1).Value = 2200
$ws.Cells.Item(25, 1).Value = 2300
$ws.Cells.Item(26, 1).Value = 2400
$ws.Cells.Item(27, 1).Value = 2500
$ws.Cells.Item(28, 1).Value = 2600
$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(30, 1).Value = 2800

# Append one new row (31) at the bottom with timestamp 2900
$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "falling"
$ws.Cells.Item(31, 3).Value = -0.143467903137207
$ws.Cells.Item(31, 4).Value = 0.759878396987915
$ws.Cells.Item(31, 5).Value = -0.6699965000152588
$ws.Cells.Item(31, 6).Value = -0.0186313893646001
$ws.Cells.Item(31, 7).Value = -0.1244637966156005
$ws.Cells.Item(31, 8).Value = 0.0003054326225537
